$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting all existing data rows down by one.
$ws.Rows.Item(2).Insert()

# The inserted row picks up formatting copied from the row above (the bold
# header row). Reset it to the plain/unstyled look used by the rest of the
# data rows, then restore the date number format on column D to match the
# other rows in that column.
$ws.Rows.Item(2).Style = "Normal"
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the new weekly record.
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C2").Value = "Coquimbo"
$ws.Range("D2").Value = 44497
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 100112022
$ws.Range("G2").Value = "Arveja Verde"
$ws.Range("H2").Value = "Perfection"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 13000
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = 14000
$ws.Range("N2").Value = '$/malla 25 kilos'
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 560
$ws.Range("Q2").Value = 25
$ws.Range("R2").Value = "Hortaliza"
